$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151, pushing the existing row 151 (and below) down to 152.
$ws.Rows.Item(151).Insert()

# Populate the new row 151 with the latest weekly entry.
$ws.Cells.Item(151, 1).Value = 1
$ws.Cells.Item(151, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(151, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(151, 4).Value = 44890
$ws.Cells.Item(151, 5).Value = 15
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100106
$ws.Cells.Item(151, 8).Value = "Oleaginosos"
$ws.Cells.Item(151, 9).Value = 100106002
$ws.Cells.Item(151, 10).Value = "Palta"
$ws.Cells.Item(151, 11).Value = "Hass"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 600
$ws.Cells.Item(151, 14).Value = 27000
$ws.Cells.Item(151, 15).Value = 28000
$ws.Cells.Item(151, 16).Value = 27500
$ws.Cells.Item(151, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(151, 18).Value = "Perú"
$ws.Cells.Item(151, 19).Value = 2750
$ws.Cells.Item(151, 20).Value = 10
